$wb = $excel.ActiveWorkbook

$newTimestamp = "2025-12-18 03:04:27"

# Sheet: Worksheets.Item(2)
$ws = $wb.Worksheets.Item(2)

# Row 10
$ws.Cells.Item(10, 3).Value = 22
$ws.Cells.Item(10, 4).Value = 382
$ws.Cells.Item(10, 5).Value = 194
$ws.Cells.Item(10, 6).Value = 188
$ws.Cells.Item(10, 7).Value = 17.36
$ws.Cells.Item(10, 8).Value = 8.82
$ws.Cells.Item(10, 9).Value = 8.550000000000001
$ws.Cells.Item(10, 10).Value = 97
$ws.Cells.Item(10, 11).Value = 84
$ws.Cells.Item(10, 23).Value = 12

# Row 11
$ws.Cells.Item(11, 3).Value = 27
$ws.Cells.Item(11, 4).Value = 622
$ws.Cells.Item(11, 5).Value = 290
$ws.Cells.Item(11, 6).Value = 332
$ws.Cells.Item(11, 7).Value = 23.04
$ws.Cells.Item(11, 8).Value = 10.74
$ws.Cells.Item(11, 9).Value = 12.3
$ws.Cells.Item(11, 10).Value = 130
$ws.Cells.Item(11, 11).Value = 121

# Row 15
$ws.Cells.Item(15, 3).Value = 22
$ws.Cells.Item(15, 4).Value = 396
$ws.Cells.Item(15, 5).Value = 189
$ws.Cells.Item(15, 6).Value = 207
$ws.Cells.Item(15, 7).Value = 18
$ws.Cells.Item(15, 8).Value = 8.59
$ws.Cells.Item(15, 9).Value = 9.41
$ws.Cells.Item(15, 10).Value = 72
$ws.Cells.Item(15, 11).Value = 91

# Row 22
$ws.Cells.Item(22, 3).Value = 23
$ws.Cells.Item(22, 4).Value = 507
$ws.Cells.Item(22, 5).Value = 242
$ws.Cells.Item(22, 6).Value = 265
$ws.Cells.Item(22, 7).Value = 22.04
$ws.Cells.Item(22, 8).Value = 10.52
$ws.Cells.Item(22, 9).Value = 11.52
$ws.Cells.Item(22, 10).Value = 91
$ws.Cells.Item(22, 11).Value = 95
$ws.Cells.Item(22, 23).Value = 6

for ($r = 2; $r -le 26; $r++) {
    $ws.Cells.Item($r, 27).Value = $newTimestamp
}

# Sheet: Worksheets.Item(3)
$ws = $wb.Worksheets.Item(3)

# Row 6
$ws.Cells.Item(6, 3).Value = 21
$ws.Cells.Item(6, 4).Value = 378
$ws.Cells.Item(6, 5).Value = 177
$ws.Cells.Item(6, 6).Value = 201
$ws.Cells.Item(6, 7).Value = 18
$ws.Cells.Item(6, 8).Value = 8.43
$ws.Cells.Item(6, 9).Value = 9.57
$ws.Cells.Item(6, 10).Value = 81
$ws.Cells.Item(6, 11).Value = 93

# Row 11
$ws.Cells.Item(11, 3).Value = 22
$ws.Cells.Item(11, 4).Value = 308
$ws.Cells.Item(11, 5).Value = 140
$ws.Cells.Item(11, 6).Value = 168
$ws.Cells.Item(11, 7).Value = 14
$ws.Cells.Item(11, 8).Value = 6.36
$ws.Cells.Item(11, 9).Value = 7.64
$ws.Cells.Item(11, 10).Value = 70
$ws.Cells.Item(11, 11).Value = 79

# Row 17
$ws.Cells.Item(17, 3).Value = 17
$ws.Cells.Item(17, 4).Value = 302
$ws.Cells.Item(17, 5).Value = 169
$ws.Cells.Item(17, 6).Value = 133
$ws.Cells.Item(17, 7).Value = 17.76
$ws.Cells.Item(17, 9).Value = 7.82
$ws.Cells.Item(17, 10).Value = 67
$ws.Cells.Item(17, 11).Value = 59
$ws.Cells.Item(17, 23).Value = 10

# Row 22
$ws.Cells.Item(22, 3).Value = 25
$ws.Cells.Item(22, 4).Value = 494
$ws.Cells.Item(22, 5).Value = 264
$ws.Cells.Item(22, 6).Value = 230
$ws.Cells.Item(22, 7).Value = 19.76
$ws.Cells.Item(22, 8).Value = 10.56
$ws.Cells.Item(22, 9).Value = 9.199999999999999
$ws.Cells.Item(22, 10).Value = 102
$ws.Cells.Item(22, 11).Value = 105
$ws.Cells.Item(22, 23).Value = 30

for ($r = 2; $r -le 26; $r++) {
    $ws.Cells.Item($r, 27).Value = $newTimestamp
}
